$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 45000
$ws.Range("J87").Value = 45000
$ws.Range("L87").Value = 45000
$ws.Range("N87").Value = -47496
$ws.Range("H90").Value = 45000
$ws.Range("J90").Value = 45000
$ws.Range("L90").Value = 135000
$ws.Range("N90").Value = -147480
$ws.Range("H137").Value = 3959
$ws.Range("I137").Value = 3449.5
$ws.Range("K137").Value = 10348.5
$ws.Range("M137").Value = -7798.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1970
$ws.Range("I102").Value = 1970
$ws.Range("K102").Value = 1970
$ws.Range("M102").Value = -348
$ws.Range("H122").Value = 3045.375
$ws.Range("I122").Value = 2894
$ws.Range("J122").Value = 3499.5
$ws.Range("K122").Value = 8682
$ws.Range("L122").Value = 10498.5
$ws.Range("M122").Value = -6232
$ws.Range("N122").Value = -15398.5
$ws.Range("H132").Value = 3491.077
$ws.Range("I132").Value = 3573.6667
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 10721.0001
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -8191.000100000001
$ws.Range("N132").Value = -12560
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 683.7368
$ws.Range("I20").Value = 453.6
$ws.Range("J20").Value = 939.44446
$ws.Range("K20").Value = 453.6
$ws.Range("L20").Value = 939.44446
$ws.Range("M20").Value = -206.6
$ws.Range("N20").Value = -1433.44446
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 550.61536
$ws.Range("I7").Value = 630.8889
$ws.Range("J7").Value = 370
$ws.Range("K7").Value = 630.8889
$ws.Range("L7").Value = 370
$ws.Range("M7").Value = -517.8889
$ws.Range("N7").Value = -596
$ws.Range("H15").Value = 3250
$ws.Range("I15").Value = 3000
$ws.Range("K15").Value = 3000
$ws.Range("M15").Value = -2830
$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H50").Value = 20270.75
$ws.Range("I50").Value = 21083
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 21083
$ws.Range("L50").Value = 20000
$ws.Range("M50").Value = -20458
$ws.Range("N50").Value = -21250
$ws.Range("H58").Value = 4722.1113
$ws.Range("I58").Value = 3500
$ws.Range("J58").Value = 5699.8
$ws.Range("K58").Value = 3500
$ws.Range("L58").Value = 5699.8
$ws.Range("M58").Value = -3297
$ws.Range("N58").Value = -6105.8
$ws.Range("H59").Value = 25318.455
$ws.Range("I59").Value = 14626
$ws.Range("J59").Value = 31428.428
$ws.Range("K59").Value = 14626
$ws.Range("L59").Value = 31428.428
$ws.Range("M59").Value = -13481
$ws.Range("N59").Value = -33718.428
$ws.Range("H68").Value = 25000
$ws.Range("I68").Value = 10000
$ws.Range("J68").Value = 30000
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 30000
$ws.Range("M68").Value = -9251
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 25000
$ws.Range("I71").Value = 10000
$ws.Range("J71").Value = 30000
$ws.Range("K71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("M71").Value = -26256
$ws.Range("N71").Value = -97488
$ws.Range("H99").Value = 4006.6155
$ws.Range("I99").Value = 3566.2222
$ws.Range("J99").Value = 4997.5
$ws.Range("K99").Value = 3566.2222
$ws.Range("L99").Value = 4997.5
$ws.Range("M99").Value = -2068.2222
$ws.Range("N99").Value = -7993.5
$ws.Range("H105").Value = 957.5714
$ws.Range("I105").Value = 1021.6
$ws.Range("K105").Value = 1021.6
$ws.Range("M105").Value = 725.4
$ws.Range("H122").Value = 2331
$ws.Range("I122").Value = 2331
$ws.Range("K122").Value = 6993
$ws.Range("M122").Value = -4543
$ws.Range("H126").Value = 4006.6155
$ws.Range("I126").Value = 3566.2222
$ws.Range("J126").Value = 4997.5
$ws.Range("K126").Value = 10698.6666
$ws.Range("L126").Value = 14992.5
$ws.Range("M126").Value = -8228.6666
$ws.Range("N126").Value = -19932.5
$ws.Range("H132").Value = 3715.0588
$ws.Range("I132").Value = 3550.5386
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 10651.6158
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -8121.6158
$ws.Range("N132").Value = -17809.25
$ws.Range("H133").Value = 29750
$ws.Range("J133").Value = 29750
$ws.Range("L133").Value = 29750
$ws.Range("N133").Value = -34810
$ws.Range("H136").Value = 4722.1113
$ws.Range("I136").Value = 3500
$ws.Range("J136").Value = 5699.8
$ws.Range("K136").Value = 10500
$ws.Range("L136").Value = 17099.4
$ws.Range("M136").Value = -7950
$ws.Range("N136").Value = -22199.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2161.7144
$ws.Range("J4").Value = 1989.8
$ws.Range("L4").Value = 5969.4
$ws.Range("N4").Value = -6193.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21640
$ws.Range("H102").Value = 307
$ws.Range("I102").Value = 307
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 307
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 1315
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 2249.5
$ws.Range("I132").Value = 2250
$ws.Range("J132").Value = 2249
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 6747
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -11807
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19821.285
$ws.Range("I7").Value = 18791.5
$ws.Range("J7").Value = 26000
$ws.Range("K7").Value = 18791.5
$ws.Range("L7").Value = 26000
$ws.Range("M7").Value = -18679.5
$ws.Range("N7").Value = -26224
$ws.Range("H122").Value = 3262.3333
$ws.Range("I122").Value = 3045.125
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9135.375
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6685.375
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 19821.285
$ws.Range("I126").Value = 18791.5
$ws.Range("J126").Value = 26000
$ws.Range("K126").Value = 56374.5
$ws.Range("L126").Value = 78000
$ws.Range("M126").Value = -53904.5
$ws.Range("N126").Value = -82940
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H122").Value = 2360.375
$ws.Range("I122").Value = 2411.8572
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7235.571599999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4785.571599999999
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 38166.332
$ws.Range("I126").Value = 31928.428
$ws.Range("K126").Value = 95785.284
$ws.Range("M126").Value = -93315.284
